$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 70, shifting existing rows 70:165 down to 71:166
$ws.Rows.Item(70).Insert()

# Populate the new row 70 with the new record's data
$ws.Range("A70").Value = 9
$ws.Range("B70").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C70").Value = "Metropolitana"
$ws.Range("D70").Value = 44967
$ws.Range("D70").NumberFormat = $ws.Range("D71").NumberFormat
$ws.Range("E70").Value = 13
$ws.Range("F70").Value = "Fruta"
$ws.Range("G70").Value = 100103
$ws.Range("H70").Value = "Frutos de hueso (carozo)"
$ws.Range("I70").Value = 100103002
$ws.Range("J70").Value = "Ciruela"
$ws.Range("K70").Value = "Larry Ann"
$ws.Range("L70").Value = "Primera"
$ws.Range("M70").Value = 200
$ws.Range("N70").Value = 10500
$ws.Range("O70").Value = 10500
$ws.Range("P70").Value = 10500
$ws.Range("Q70").Value = "$/caja 15 kilos granel"
$ws.Range("R70").Value = "Región de O'Higgins"
$ws.Range("S70").Value = 700
$ws.Range("T70").Value = 15
